$wb = $excel.ActiveWorkbook

$ws3 = $wb.Worksheets.Item("RO & CO Hearing Allocation")

$ws3.Range("A1").Value = "Allocation of Regional Office Video Hearings"

$ws3.Rows.Item(4).Delete()

$ws3.Name = "RO Allocations"
